$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

#############################################
# 1) Header text edits (surgical in shared strings)
#############################################

# A8: "Volume 31   Number  13" -> "...14" (Volume/Number text)
$a8 = $ws.Range("A8")
$a8full = $a8.Value()
$a8idx = $a8full.LastIndexOf("13")
$a8.Characters($a8idx + 1, 2).Text = "14"

# C9: "Report Covering the Week  3/25/2024  Through  3/31/2024" -> "...4/1/2024  Through  4/7/2024"
$c9 = $ws.Range("C9")
$c9full = $c9.Value()
$c9idx1 = $c9full.IndexOf("3/25/2024")
$c9.Characters($c9idx1 + 1, 9).Text = "4/1/2024"
$c9full2 = $ws.Range("C9").Value()
$c9idx2 = $c9full2.IndexOf("3/31/2024")
$c9.Characters($c9idx2 + 1, 9).Text = "4/7/2024"

#############################################
# 2) Cells that switch between numeric <-> shared-string "blank" (0 / ***.*)
#    Use Copy() from a same-style donor cell in row 14 to carry over the
#    correct style index + shared-string type without creating new styles,
#    then overwrite the value where a literal number is needed.
#############################################

$ws.Range("C14").Copy($ws.Range("F15"))

$ws.Range("G14").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 5

$ws.Range("C14").Copy($ws.Range("C22"))

$ws.Range("C14").Copy($ws.Range("F27"))

$ws.Range("G14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1

$ws.Range("C14").Copy($ws.Range("D28"))

$ws.Range("E14").Copy($ws.Range("E28"))

#############################################
# 3) Plain numeric value updates (style/type unchanged)
#############################################

$ws.Range("N15").Value = -44.444444444444
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -7.142857142857
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 40
$ws.Range("L16").Value = -9.259259259259
$ws.Range("M16").Value = -33.783783783783
$ws.Range("N16").Value = -85.060975609756
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -36.363636363636
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 103
$ws.Range("K17").Value = -24.271844660194
$ws.Range("L17").Value = -6.024096385542
$ws.Range("M17").Value = 36.842105263157
$ws.Range("N17").Value = -8.235294117647
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = -28.571428571428
$ws.Range("L18").Value = -19.354838709677
$ws.Range("M18").Value = -64.788732394366
$ws.Range("N18").Value = -92.957746478873
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 5.263157894736
$ws.Range("I19").Value = 93
$ws.Range("J19").Value = 84
$ws.Range("K19").Value = 10.714285714285
$ws.Range("L19").Value = -13.084112149532
$ws.Range("M19").Value = 16.25
$ws.Range("N19").Value = -30.597014925373
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 59
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = -16.901408450704
$ws.Range("L20").Value = -4.838709677419
$ws.Range("M20").Value = -10.606060606060
$ws.Range("N20").Value = -93.904958677686
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -22.222222222222
$ws.Range("F21").Value = 81
$ws.Range("H21").Value = -13.829787234042
$ws.Range("I21").Value = 309
$ws.Range("J21").Value = 333
$ws.Range("K21").Value = -7.207207207207
$ws.Range("L21").Value = -9.117647058823
$ws.Range("M21").Value = -11.461318051575
$ws.Range("N21").Value = -83.581296493092
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 87.5
$ws.Range("F24").Value = 116
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = 1.754385964912
$ws.Range("I24").Value = 375
$ws.Range("J24").Value = 397
$ws.Range("K24").Value = -5.541561712846
$ws.Range("L24").Value = 5.042016806722
$ws.Range("M24").Value = 80.288461538461
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 80.645161290322
$ws.Range("I25").Value = 179
$ws.Range("J25").Value = 126
$ws.Range("K25").Value = 42.063492063492
$ws.Range("L25").Value = 29.710144927536
$ws.Range("C26").Value = 15
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 55
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 161
$ws.Range("J26").Value = 143
$ws.Range("K26").Value = 12.587412587412
$ws.Range("L26").Value = 27.777777777777
$ws.Range("M26").Value = 8.053691275167
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = -12.5
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("I28").Value = 10
$ws.Range("K28").Value = -37.5
$ws.Range("L28").Value = -23.076923076923
$ws.Range("G29").Value = 3
$ws.Range("G30").Value = 3
$ws.Range("L31").Value = 0
